# Horarios Linea 141 - actualizacion scrap 15:16:46 (Total filas sheet1: 212 -> 219)
# Logic: each sheet holds scraped rows (Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada).
# A fresh scrape appends new rows (all sharing the new Hora_Scrap timestamp) and the whole
# data block gets re-sorted (stable) by column B (Hora_Llegada) ascending - mirrors the
# scraper's "re-dump sorted by arrival time" behaviour visible in the diff.

$wb = $excel.ActiveWorkbook

function TimeToSeconds($t) {
    $parts = $t -split ":"
    $h = [int]$parts[0]
    $m = [int]$parts[1]
    $s = 0
    if ($parts.Length -gt 2) { $s = [int]$parts[2] }
    return ($h * 3600) + ($m * 60) + $s
}

function ReadRows($ws, $firstRow, $lastRow) {
    $rows = @()
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $a = $ws.Cells.Item($r, 1).Value2
        $b = $ws.Cells.Item($r, 2).Value2
        $c = $ws.Cells.Item($r, 3).Value2
        $d = $ws.Cells.Item($r, 4).Value2
        $e = $ws.Cells.Item($r, 5).Value2
        $rows += [PSCustomObject]@{ A = $a; B = $b; C = $c; D = $d; E = $e; SortKey = (TimeToSeconds $b) }
    }
    return $rows
}

function WriteRows($ws, $firstRow, $rows) {
    $r = $firstRow
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row.A
        $ws.Cells.Item($r, 2).Value = $row.B
        $ws.Cells.Item($r, 3).Value = $row.C
        $ws.Cells.Item($r, 4).Value = $row.D
        $ws.Cells.Item($r, 5).Value = $row.E
        $r = $r + 1
    }
}

function NewRow($a, $b, $c, $d, $e) {
    return [PSCustomObject]@{ A = $a; B = $b; C = $c; D = $d; E = $e; SortKey = (TimeToSeconds $b) }
}

$scrapTime = "15:16:46"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$existing1 = ReadRows $ws1 6 217

$newRows1 = @(
    (NewRow $scrapTime "16:05" "16_SANTA ANA" 49 "LP1912"),
    (NewRow $scrapTime "16:22" "23_HERNANDEZ" 66 "LP1912"),
    (NewRow $scrapTime "16:30" "16_SANTA ANA" 74 "LP1912"),
    (NewRow $scrapTime "16:53" "11_ETCHEVERRY" 97 "LP1912"),
    (NewRow $scrapTime "16:58" "15_ABASTO" 102 "LP1912"),
    (NewRow $scrapTime "17:07" "16_P MOR-SANTA ANA" 111 "LP1912"),
    (NewRow $scrapTime "17:09" "215C_EL PATO" 113 "LP1912")
)

$combined1 = @($existing1) + @($newRows1)
$sorted1 = $combined1 | Sort-Object -Property SortKey

WriteRows $ws1 6 $sorted1

$total1 = $sorted1.Length
$lastRow1 = 5 + $total1

$ws1.Range("A2").Value = "Última actualización: $scrapTime"
$ws1.Range("A3").Value = "Total filas: $total1"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215 (subset of sheet1 rows whose Linea mentions "215")
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$sheet2Rows = @()
foreach ($row in $sorted1) {
    if ($row.C -like "*215*") {
        $sheet2Rows += (NewRow $row.A $row.B $row.C $row.D $row.E)
    }
}

WriteRows $ws2 6 $sheet2Rows

$total2 = $sheet2Rows.Length
$ws2.Range("A2").Value = "Última actualización: $scrapTime"
$ws2.Range("A3").Value = "Total filas: $total2"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$existing3 = ReadRows $ws3 6 40

$newRows3 = @(
    (NewRow $scrapTime "15:19" "215A_LA PLATA" 3 "L6173"),
    (NewRow $scrapTime "17:05" "215C_LA PLATA" 109 "L6203")
)

$combined3 = @($existing3) + @($newRows3)
$sorted3 = $combined3 | Sort-Object -Property SortKey

WriteRows $ws3 6 $sorted3

$total3 = $sorted3.Length
$ws3.Range("A2").Value = "Última actualización: $scrapTime"
$ws3.Range("A3").Value = "Total filas: $total3"

Write-Host "sheet1 rows: $total1, sheet2 rows: $total2, sheet3 rows: $total3"
